$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAY_2024")

# Update attendance values for the newly attended class (column E / D)
$ws.Range("E14").Value = 3
$ws.Range("E15").Value = 3
$ws.Range("E16").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E18").Value = 3
$ws.Range("D21").Value = 3

# Force recalculation of dependent formulas (G and H columns)
$excel.Calculate()

# Update the active cell selection to match the saved view state
$ws.Activate()
$ws.Range("E19").Select()
